$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the used range with the same base cell style (s="1") that the
# --- existing "Roll number" column already uses, for every new cell we are
# --- about to populate (new "Amount" column B1:B10, and the two extra
# --- duplicate roll-number rows A7:A8, plus the new trailing blank rows).
$ws.Range("A1").Copy()
$ws.Range("B1:B10").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A7:A10").PasteSpecial(-4122)

# --- New "Amount" header + values next to the existing "Roll number" data
$ws.Range("B1").Value = "Amount"
$ws.Range("B2").Value = 11
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 13
$ws.Range("B5").Value = 15
$ws.Range("B6").Value = 25

# --- New duplicate roll-number entries ("20bcs001" appears twice) with dues
$ws.Range("A7").Value = "20bcs001"
$ws.Range("B7").Value = 3
$ws.Range("A8").Value = "20bcs001"
$ws.Range("B8").Value = 9

# --- Highlight the duplicate's amount cell: white fill, black Arial font,
# --- right aligned - flags the duplicate entry for integrity review.
$ws.Range("B8").Font.Color = 0
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B8").Interior.Color = 16777215
$ws.Range("B8").HorizontalAlignment = -4152

# --- Rows 9 & 10 stay blank (placeholders reserved for future dues entries)
# --- but already carry formatting from the PasteSpecial above, so nothing
# --- else to set here.
